$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-17"

# Update the label for the March row
$ws.Range("A4").Value = "March (through 03-17)"

# Update March row (row 4) values
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 35
$ws.Range("F4").Value = 17
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 74

# Update Total row (row 5) values
$ws.Range("B5").Value = 54
$ws.Range("C5").Value = 108
$ws.Range("D5").Value = 163
$ws.Range("E5").Value = 172
$ws.Range("F5").Value = 96
$ws.Range("G5").Value = 172
$ws.Range("H5").Value = 389
$ws.Range("I5").Value = 374
